$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture B4's current ("date only") number format before it gets overwritten
# below, so the newly appended row can reuse it.
$dateOnlyFormat = $ws.Range("B4").NumberFormat

# Row 4, column B (date) switches to the "date+time" format used by rows 2/3.
$ws.Range("B4").NumberFormat = $ws.Range("B2").NumberFormat

# Append a new row 5 with values, reusing the "date only" format that row 4's
# B cell previously had.
$ws.Range("A5").Value = 803.816
$ws.Range("B5").Value = 45729
$ws.Range("B5").NumberFormat = $dateOnlyFormat
$ws.Range("C5").Value = 810.465
